$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).Value = "2025-08-21 06:56:24 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-21 12:26:24 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

# Match the formatting of the previous data row (style s="3": centered, same font/border)
$ws.Range("A38:H38").Copy()
$ws.Range("A39:H39").PasteSpecial(-4122)
